$wb = $excel.ActiveWorkbook

# ============ Sheet 1 ============
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "admin1"
$ws.Range("B1").Value = "pop_group"

# Column A (admin1 codes)
$ws.Range("A2").Value = "MMR001"
$ws.Range("A3").Value = "MMR002"
$ws.Range("A4").Value = "MMR003"
$ws.Range("A5").Value = "MMR004"
$ws.Range("A6").Value = "MMR005"
$ws.Range("A7").Value = "MMR006"
$ws.Range("A8").Value = "MMR007"
$ws.Range("A9").Value = "MMR009"
$ws.Range("A10").Value = "MMR011"
$ws.Range("A11").Value = "MMR012"
$ws.Range("A12").Value = "MMR014"
$ws.Range("A13").Value = "MMR015"

# Column B (pop_group)
$ws.Range("B2").Value = "idp"
$ws.Range("B3").Value = "idp"
$ws.Range("B4").Value = "idp"
$ws.Range("B5").Value = "idp"
$ws.Range("B6").Value = "idp"
$ws.Range("B7").Value = "idp"
$ws.Range("B8").Value = "idp"
$ws.Range("B9").Value = "idp"
$ws.Range("B10").Value = "idp"
$ws.Range("B11").Value = "idp"
$ws.Range("B12").Value = "idp"
$ws.Range("B13").Value = "idp"

# Numeric columns C-F
$ws.Range("C2").Value = 0.6540118601949976
$ws.Range("D2").Value = 0.2329490092862201
$ws.Range("E2").Value = 0.1130391305187824
$ws.Range("F2").Value = 0
$ws.Range("C3").Value = 0.3195841598793305
$ws.Range("D3").Value = 0.5897774930562981
$ws.Range("E3").Value = 0.02163833014917246
$ws.Range("F3").Value = 0.06900001691519914
$ws.Range("C4").Value = 0.4282774371907243
$ws.Range("D4").Value = 0.465292329157762
$ws.Range("E4").Value = 0.05867773143280703
$ws.Range("F4").Value = 0.04775250221870657
$ws.Range("C5").Value = 0.2596001116751376
$ws.Range("D5").Value = 0.7251525007435389
$ws.Range("E5").Value = 0.01524738758132347
$ws.Range("F5").Value = 0
$ws.Range("C6").Value = 0.04861766279749088
$ws.Range("D6").Value = 0.8488057624659348
$ws.Range("E6").Value = 0.01874899599464879
$ws.Range("F6").Value = 0.08382757874192558
$ws.Range("C7").Value = 0.3139748305934238
$ws.Range("D7").Value = 0.5915523224995896
$ws.Range("E7").Value = 0.06789595527463371
$ws.Range("F7").Value = 0.0265768916323528
$ws.Range("C8").Value = 0.6688827647836127
$ws.Range("D8").Value = 0.2618557855553887
$ws.Range("E8").Value = 0.06926144966099851
$ws.Range("F8").Value = 0
$ws.Range("C9").Value = 0.4228032570458399
$ws.Range("D9").Value = 0.4106630319422648
$ws.Range("E9").Value = 0.08301522853837168
$ws.Range("F9").Value = 0.08351848247352356
$ws.Range("C10").Value = 0.5553267378807409
$ws.Range("D10").Value = 0.2857939673992496
$ws.Range("E10").Value = 0.07349868043284102
$ws.Range("F10").Value = 0.08538061428716844
$ws.Range("C11").Value = 0.2903546969040392
$ws.Range("D11").Value = 0.5095981950530968
$ws.Range("E11").Value = 0.1767228389255024
$ws.Range("F11").Value = 0.02332426911736172
$ws.Range("C12").Value = 0.5593708819099853
$ws.Range("D12").Value = 0.4086532998696034
$ws.Range("E12").Value = 0.03197581822041135
$ws.Range("F12").Value = 0
$ws.Range("C13").Value = 0.1583284297571598
$ws.Range("D13").Value = 0.8068898558470062
$ws.Range("E13").Value = 0.03478171439583397
$ws.Range("F13").Value = 0

# ============ Sheet 2 ============
$ws = $wb.Worksheets.Item(2)
$ws.Range("A3:F6").ClearContents()
$ws.Range("A1").Value = "admin1"
$ws.Range("B1").Value = "pop_group"

# Column A (admin1 codes)
$ws.Range("A2").Value = "MMR012"

# Column B (pop_group)
$ws.Range("B2").Value = "ndsp"

# Numeric columns C-F
$ws.Range("C2").Value = 0.2823734803505052
$ws.Range("D2").Value = 0.6573722163971507
$ws.Range("E2").Value = 0.06025430325234421
$ws.Range("F2").Value = 0

# ============ Sheet 3 ============
$ws = $wb.Worksheets.Item(3)
$ws.Range("A1").Value = "admin1"
$ws.Range("B1").Value = "pop_group"

# Column A (admin1 codes)
$ws.Range("A2").Value = "MMR001"
$ws.Range("A3").Value = "MMR002"
$ws.Range("A4").Value = "MMR003"
$ws.Range("A5").Value = "MMR004"
$ws.Range("A6").Value = "MMR005"
$ws.Range("A7").Value = "MMR006"
$ws.Range("A8").Value = "MMR007"
$ws.Range("A9").Value = "MMR008"
$ws.Range("A10").Value = "MMR009"
$ws.Range("A11").Value = "MMR010"
$ws.Range("A12").Value = "MMR011"
$ws.Range("A13").Value = "MMR012"
$ws.Range("A14").Value = "MMR013"
$ws.Range("A15").Value = "MMR014"
$ws.Range("A16").Value = "MMR015"
$ws.Range("A17").Value = "MMR016"
$ws.Range("A18").Value = "MMR017"
$ws.Range("A19").Value = "MMR018"

# Column B (pop_group)
$ws.Range("B2").Value = "ocap"
$ws.Range("B3").Value = "ocap"
$ws.Range("B4").Value = "ocap"
$ws.Range("B5").Value = "ocap"
$ws.Range("B6").Value = "ocap"
$ws.Range("B7").Value = "ocap"
$ws.Range("B8").Value = "ocap"
$ws.Range("B9").Value = "ocap"
$ws.Range("B10").Value = "ocap"
$ws.Range("B11").Value = "ocap"
$ws.Range("B12").Value = "ocap"
$ws.Range("B13").Value = "ocap"
$ws.Range("B14").Value = "ocap"
$ws.Range("B15").Value = "ocap"
$ws.Range("B16").Value = "ocap"
$ws.Range("B17").Value = "ocap"
$ws.Range("B18").Value = "ocap"
$ws.Range("B19").Value = "ocap"

# Numeric columns C-F
$ws.Range("C2").Value = 0.7479750560573946
$ws.Range("D2").Value = 0.229998866780615
$ws.Range("E2").Value = 0.02202607716199041
$ws.Range("F2").Value = 0
$ws.Range("C3").Value = 0.2899498531625629
$ws.Range("D3").Value = 0.6624921725493361
$ws.Range("E3").Value = 0.0475579742881008
$ws.Range("F3").Value = 0
$ws.Range("C4").Value = 0.5664132847520789
$ws.Range("D4").Value = 0.3802086730760768
$ws.Range("E4").Value = 0.04346546589109143
$ws.Range("F4").Value = 0.009912576280752863
$ws.Range("C5").Value = 0.2763671132366832
$ws.Range("D5").Value = 0.6550785339776489
$ws.Range("E5").Value = 0.05162703592585307
$ws.Range("F5").Value = 0.01692731685981476
$ws.Range("C6").Value = 0.1746538347093081
$ws.Range("D6").Value = 0.773964089295205
$ws.Range("E6").Value = 0.04132606974079942
$ws.Range("F6").Value = 0.01005600625468734
$ws.Range("C7").Value = 0.7936303282700614
$ws.Range("D7").Value = 0.1870216821339076
$ws.Range("E7").Value = 0.01934798959603097
$ws.Range("F7").Value = 0
$ws.Range("C8").Value = 0.7312961625531302
$ws.Range("D8").Value = 0.2349945367094737
$ws.Range("E8").Value = 0.03370930073739602
$ws.Range("F8").Value = 0
$ws.Range("C9").Value = 0.7801742695693389
$ws.Range("D9").Value = 0.1575988297589671
$ws.Range("E9").Value = 0.0622269006716941
$ws.Range("F9").Value = 0
$ws.Range("C10").Value = 0.7929942800969818
$ws.Range("D10").Value = 0.1780685290760373
$ws.Range("E10").Value = 0.02893719082698087
$ws.Range("F10").Value = 0
$ws.Range("C11").Value = 0.791065957409446
$ws.Range("D11").Value = 0.1931232883155541
$ws.Range("E11").Value = 0.01581075427499976
$ws.Range("F11").Value = 0
$ws.Range("C12").Value = 0.628663794284451
$ws.Range("D12").Value = 0.1764936784966122
$ws.Range("E12").Value = 0.09349543837260667
$ws.Range("F12").Value = 0.1013470888463301
$ws.Range("C13").Value = 0.4391019149385846
$ws.Range("D13").Value = 0.4876446382271574
$ws.Range("E13").Value = 0.06455350389900151
$ws.Range("F13").Value = 0.008699942935256486
$ws.Range("C14").Value = 0.624700271132733
$ws.Range("D14").Value = 0.2918746545088312
$ws.Range("E14").Value = 0.07539968673165817
$ws.Range("F14").Value = 0.00802538762677754
$ws.Range("C15").Value = 0.7841582527586344
$ws.Range("D15").Value = 0.1496428406259407
$ws.Range("E15").Value = 0.06619890661542482
$ws.Range("F15").Value = 0
$ws.Range("C16").Value = 0.5831172771890769
$ws.Range("D16").Value = 0.3873957991753649
$ws.Range("E16").Value = 0.02948692363555816
$ws.Range("F16").Value = 0
$ws.Range("C17").Value = 0.8444730632458506
$ws.Range("D17").Value = 0.1165588997972592
$ws.Range("E17").Value = 0.03896803695689039
$ws.Range("F17").Value = 0
$ws.Range("C18").Value = 0.9240234766960284
$ws.Range("D18").Value = 0.026949994451861
$ws.Range("E18").Value = 0.04902652885211051
$ws.Range("F18").Value = 0
$ws.Range("C19").Value = 0.8674312585554411
$ws.Range("D19").Value = 0.1053186169796036
$ws.Range("E19").Value = 0.02725012446495541
$ws.Range("F19").Value = 0

# ============ Sheet 4 ============
$ws = $wb.Worksheets.Item(4)
$ws.Range("A14:F14").ClearContents()
$ws.Range("A1").Value = "admin1"
$ws.Range("B1").Value = "pop_group"

# Column A (admin1 codes)
$ws.Range("A2").Value = "MMR001"
$ws.Range("A3").Value = "MMR002"
$ws.Range("A4").Value = "MMR003"
$ws.Range("A5").Value = "MMR004"
$ws.Range("A6").Value = "MMR005"
$ws.Range("A7").Value = "MMR006"
$ws.Range("A8").Value = "MMR007"
$ws.Range("A9").Value = "MMR009"
$ws.Range("A10").Value = "MMR011"
$ws.Range("A11").Value = "MMR012"
$ws.Range("A12").Value = "MMR014"
$ws.Range("A13").Value = "MMR015"

# Column B (pop_group)
$ws.Range("B2").Value = "ret"
$ws.Range("B3").Value = "ret"
$ws.Range("B4").Value = "ret"
$ws.Range("B5").Value = "ret"
$ws.Range("B6").Value = "ret"
$ws.Range("B7").Value = "ret"
$ws.Range("B8").Value = "ret"
$ws.Range("B9").Value = "ret"
$ws.Range("B10").Value = "ret"
$ws.Range("B11").Value = "ret"
$ws.Range("B12").Value = "ret"
$ws.Range("B13").Value = "ret"

# Numeric columns C-F
$ws.Range("C2").Value = 0.5853868264010187
$ws.Range("D2").Value = 0.3794977558531176
$ws.Range("E2").Value = 0.03511541774586386
$ws.Range("F2").Value = 0
$ws.Range("C3").Value = 0.1881113598201735
$ws.Range("D3").Value = 0.6803408258131957
$ws.Range("E3").Value = 0.02539776024089409
$ws.Range("F3").Value = 0.1061500541257369
$ws.Range("C4").Value = 0.4522314307539207
$ws.Range("D4").Value = 0.535172232157055
$ws.Range("E4").Value = 0.01259633708902428
$ws.Range("F4").Value = 0
$ws.Range("C5").Value = 0.3628408843989914
$ws.Range("D5").Value = 0.5612727249817211
$ws.Range("E5").Value = 0.07588639061928751
$ws.Range("F5").Value = 0
$ws.Range("C6").Value = 0.07316629031944184
$ws.Range("D6").Value = 0.9268337096805581
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("C7").Value = 0.443996661956498
$ws.Range("D7").Value = 0.4760960270901249
$ws.Range("E7").Value = 0.07990731095337696
$ws.Range("F7").Value = 0
$ws.Range("C8").Value = 0.3918059515410319
$ws.Range("D8").Value = 0.534593454033195
$ws.Range("E8").Value = 0.07360059442577312
$ws.Range("F8").Value = 0
$ws.Range("C9").Value = 0.3246193244816443
$ws.Range("D9").Value = 0.5008968621725338
$ws.Range("E9").Value = 0.129743185531197
$ws.Range("F9").Value = 0.04474062781462478
$ws.Range("C10").Value = 0.4928355615432358
$ws.Range("D10").Value = 0.3626867371558554
$ws.Range("E10").Value = 0.1190506602479244
$ws.Range("F10").Value = 0.02542704105298452
$ws.Range("C11").Value = 0.2535230518171918
$ws.Range("D11").Value = 0.6450779162767399
$ws.Range("E11").Value = 0.09519507690366545
$ws.Range("F11").Value = 0.006203955002402782
$ws.Range("C12").Value = 0.5718279007459466
$ws.Range("D12").Value = 0.2237822709111292
$ws.Range("E12").Value = 0.2043898283429242
$ws.Range("F12").Value = 0
$ws.Range("C13").Value = 0.4177383882796832
$ws.Range("D13").Value = 0.4557602641860585
$ws.Range("E13").Value = 0.03652022065847701
$ws.Range("F13").Value = 0.08998112687578116
